# Feat: add roles table
#
# Adds a new "Rol" table (id, rol) to the relational-model sheet and turns
# the existing Usuario.rol column (a free-text CK enum) into a foreign key
# (idRol -> Rol.id).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New "Rol" table at J27:K30 --------------------------------------
# Clone formatting from the same-shaped "Contenedor" table block (title in
# column B, 2-row header/type block, blank separator row) onto J:K.
$ws.Range("B27").Copy()
$ws.Range("J27").PasteSpecial(-4122)
$ws.Range("B28:C30").Copy()
$ws.Range("J28:K30").PasteSpecial(-4122)

# Table title
$ws.Range("J27").Value = "Rol"

# --- 2. Usuario.rol -> Usuario.idRol (FK to Rol.id) ---------------------
$ws.Range("H24").Value = "NN, FK Rol.id"
$ws.Range("H23").Value = "idRol"

# --- 3. Fill in the rest of the new Rol table ---------------------------
$ws.Range("K29").Value = "ND, NN"
$ws.Range("J28").Value = "id"
$ws.Range("K28").Value = "rol"
$ws.Range("J29").Value = "PK, SA"

# --- 4. Row-height tweaks that followed the content changes -------------
$ws.Rows.Item(24).RowHeight = 17
$ws.Rows.Item(29).RowHeight = 72

# --- 5. Selection / scroll position left by the edit ---------------------
$ws.Range("K31").Select()
